# Applies the cryptos.xlsx price/volume refresh described in the commit
# 'Updated cryptos list on Sun May 14 08:57:42 UTC 2023 with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, volume % text) - safe to assign directly.
$textUpdates = @(
    @{ Cell = 'D2'; Value = '27.779.82' }
    @{ Cell = 'E2'; Value = '  +2.67%  ' }
    @{ Cell = 'D3'; Value = '1.863.76' }
    @{ Cell = 'E3'; Value = '  +2.27%  ' }
    @{ Cell = 'E4'; Value = '  +2.38%  ' }
    @{ Cell = 'E5'; Value = '  +3.21%  ' }
    @{ Cell = 'E6'; Value = '  +2.42%  ' }
    @{ Cell = 'E7'; Value = '  +2.25%  ' }
    @{ Cell = 'E8'; Value = '  +2.33%  ' }
    @{ Cell = 'E9'; Value = '  +2.38%  ' }
    @{ Cell = 'E10'; Value = '  +1.32%  ' }
    @{ Cell = 'E11'; Value = '  +1.64%  ' }
    @{ Cell = 'D12'; Value = '1.890.06' }
    @{ Cell = 'E12'; Value = '  -12.26%  ' }
    @{ Cell = 'E13'; Value = '  +2.34%  ' }
    @{ Cell = 'E14'; Value = '  +1.37%  ' }
    @{ Cell = 'E15'; Value = '  +2.84%  ' }
    @{ Cell = 'E16'; Value = '  +2.99%  ' }
    @{ Cell = 'E17'; Value = '  +2.25%  ' }
    @{ Cell = 'E18'; Value = '  +2.12%  ' }
    @{ Cell = 'E19'; Value = '  +2.58%  ' }
    @{ Cell = 'E20'; Value = '  +1.64%  ' }
    @{ Cell = 'D21'; Value = '27.794.37' }
    @{ Cell = 'E21'; Value = '  +2.52%  ' }
    @{ Cell = 'E22'; Value = '  +1.99%  ' }
    @{ Cell = 'E23'; Value = '  +2.97%  ' }
    @{ Cell = 'E24'; Value = '  +3.67%  ' }
    @{ Cell = 'E25'; Value = '  +2.32%  ' }
    @{ Cell = 'E26'; Value = '  +2.21%  ' }
    @{ Cell = 'E27'; Value = '  +3.31%  ' }
    @{ Cell = 'E28'; Value = '  +1.31%  ' }
    @{ Cell = 'E29'; Value = '  +2.07%  ' }
    @{ Cell = 'E30'; Value = '  +1.53%  ' }
    @{ Cell = 'B31'; Value = 'ImmutableX' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'E31'; Value = '  +3.55%  ' }
    @{ Cell = 'B32'; Value = 'ARBITRUM' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = 'E32'; Value = '  +3.40%  ' }
    @{ Cell = 'E33'; Value = '  +8.83%  ' }
    @{ Cell = 'E34'; Value = '  +2.98%  ' }
    @{ Cell = 'E35'; Value = '  +2.68%  ' }
    @{ Cell = 'E36'; Value = '  +3.48%  ' }
    @{ Cell = 'E37'; Value = '  +3.23%  ' }
    @{ Cell = 'E38'; Value = '  +1.96%  ' }
    @{ Cell = 'B39'; Value = 'MXToken' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'E39'; Value = '  +3.19%  ' }
    @{ Cell = 'B40'; Value = 'TheSandbox' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' }
    @{ Cell = 'E40'; Value = '  +1.09%  ' }
    @{ Cell = 'E41'; Value = '  +2.35%  ' }
    @{ Cell = 'E42'; Value = '  +6.35%  ' }
    @{ Cell = 'E43'; Value = '  +4.29%  ' }
    @{ Cell = 'E44'; Value = '  +2.52%  ' }
    @{ Cell = 'E45'; Value = '  +2.03%  ' }
    @{ Cell = 'E46'; Value = '  +3.81%  ' }
    @{ Cell = 'E47'; Value = '  +2.37%  ' }
    @{ Cell = 'E48'; Value = '  +1.87%  ' }
    @{ Cell = 'E49'; Value = '  +3.11%  ' }
    @{ Cell = 'E50'; Value = '  +4.49%  ' }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Price updates that look numeric (e.g. '1.038', '83.71') must stay text cells,
# exactly like the rest of column D, so force the Text number format while writing
# them and then drop back to the sheet's normal (unstyled) look.
$numericLookingUpdates = @(
    @{ Cell = 'D4'; Value = '1.038' }
    @{ Cell = 'D5'; Value = '324.30' }
    @{ Cell = 'D7'; Value = '0.4414' }
    @{ Cell = 'D9'; Value = '0.07461' }
    @{ Cell = 'D10'; Value = '0.8841' }
    @{ Cell = 'D13'; Value = '5.556' }
    @{ Cell = 'D14'; Value = '6.751' }
    @{ Cell = 'D15'; Value = '0.07202' }
    @{ Cell = 'D16'; Value = '83.71' }
    @{ Cell = 'D18'; Value = '0.000009138' }
    @{ Cell = 'D19'; Value = '1.035' }
    @{ Cell = 'D20'; Value = '15.54' }
    @{ Cell = 'D22'; Value = '5.318' }
    @{ Cell = 'D25'; Value = '158.43' }
    @{ Cell = 'D26'; Value = '18.88' }
    @{ Cell = 'D27'; Value = '1.992' }
    @{ Cell = 'D28'; Value = '5.322' }
    @{ Cell = 'D29'; Value = '117.59' }
    @{ Cell = 'D30'; Value = '0.09106' }
    @{ Cell = 'D31'; Value = '0.7773' }
    @{ Cell = 'D32'; Value = '1.219' }
    @{ Cell = 'D33'; Value = '3.074' }
    @{ Cell = 'D34'; Value = '4.586' }
    @{ Cell = 'D35'; Value = '1.037' }
    @{ Cell = 'D36'; Value = '1.166' }
    @{ Cell = 'D37'; Value = '0.01995' }
    @{ Cell = 'D38'; Value = '0.05358' }
    @{ Cell = 'D39'; Value = '2.848' }
    @{ Cell = 'D40'; Value = '0.5205' }
    @{ Cell = 'D41'; Value = '0.1697' }
    @{ Cell = 'D42'; Value = '6.922' }
    @{ Cell = 'D43'; Value = '8.721' }
    @{ Cell = 'D45'; Value = '109.68' }
    @{ Cell = 'D46'; Value = '1.723' }
    @{ Cell = 'D48'; Value = '0.06441' }
    @{ Cell = 'D49'; Value = '1.883' }
    @{ Cell = 'D50'; Value = '39.98' }
    @{ Cell = 'D51'; Value = '64.58' }
)

foreach ($u in $numericLookingUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
